$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crawl timestamp for every scraped row (O2:O527)
$ws.Range("O2:O527").Value = "2023-01-25 12:58:24"

# Rating-amount ("ratingAmount", column D) upticks picked up by the crawler
$ws.Range("D11").Value = 272
$ws.Range("D47").Value = 100
$ws.Range("D65").Value = 16
$ws.Range("D137").Value = 19
$ws.Range("D146").Value = 5
$ws.Range("D319").Value = 38
$ws.Range("D457").Value = 6
$ws.Range("D488").Value = 18

# Row 485 also had its ratingValue adjusted
$ws.Range("D485").Value = 5
$ws.Range("E485").Value = 3.5

# Row 383: product now shown as out of stock online
$ws.Range("M383").Value = "Naturaplan Bio Shiitake-Pilze ca. 100g - Online kein Bestand 3.20 Schweizer Franken"
